$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> D, L, M, N, O, P, Q, R, S, T values (target state per diff)
$data = @{
  2  = @(44204, "Primera", 110, 7000,  7500,  7318,  '$/bandeja 7 kilos', "Provincia de San Felipe de Aconcagua", 1045, 7)
  3  = @(44189, "Especial", 20, 15000, 15000, 15000, '$/bandeja 7 kilos', "Provincia de San Felipe de Aconcagua", 2143, 7)
  4  = @(44189, "Primera", 30, 13000, 13000, 13000, '$/bandeja 7 kilos', "Provincia de San Felipe de Aconcagua", 1857, 7)
  5  = @(44550, "Primera", 60, 24000, 24000, 24000, '$/bandeja 7 kilos', "Región Metropolitana", 3429, 7)
  6  = @(44553, "Especial", 200, 22000, 22000, 22000, '$/bandeja 6 kilos', "Provincia de San Felipe de Aconcagua", 3667, 6)
  7  = @(44553, "Primera", 150, 18000, 18000, 18000, '$/bandeja 6 kilos', "Provincia de San Felipe de Aconcagua", 3000, 6)
  8  = @(44572, "Primera", 65, 20000, 20000, 20000, '$/bandeja 6 kilos', "Región Metropolitana", 3333, 6)
  11 = @(44561, "Primera", 200, 18000, 18000, 18000, '$/bandeja 6 kilos', "Provincia de San Felipe de Aconcagua", 3000, 6)
  12 = @(44558, "Especial", 20, 22000, 22000, 22000, '$/bandeja 6 kilos', "Provincia de San Felipe de Aconcagua", 3667, 6)
  13 = @(44558, "Primera", 25, 18000, 18000, 18000, '$/bandeja 6 kilos', "Provincia de San Felipe de Aconcagua", 3000, 6)
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $ws.Range("D$row").Value = $vals[0]
  $ws.Range("L$row").Value = $vals[1]
  $ws.Range("M$row").Value = $vals[2]
  $ws.Range("N$row").Value = $vals[3]
  $ws.Range("O$row").Value = $vals[4]
  $ws.Range("P$row").Value = $vals[5]
  $ws.Range("Q$row").Value = $vals[6]
  $ws.Range("R$row").Value = $vals[7]
  $ws.Range("S$row").Value = $vals[8]
  $ws.Range("T$row").Value = $vals[9]
}
